$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trivia rows to append below the existing data (rows 58-70).
$rows = @(
    @("6909 Languages IN THE World", "Geography"),
    @("3 Periods IN A Hockey Game", "Sports"),
    @("4 Quarters IN A Basketball Game", "Sports"),
    @("4 Periods IN A Lacrosse Game", "Sports"),
    @("5 Positions IN Basketball", "Sports"),
    @("4840 Square Yards IN A Acre", "Measurements"),
    @("10000 Square Meters IN A Hectare", "Measurements"),
    @("4 Grand Slams IN Tennis", "Sports"),
    @("8 Stones IN Curling", "Sports"),
    @("15 Sports IN THE Winter Olympics", "Sports"),
    @("28 Sports IN THE Summer Olympics", "Sports"),
    @("23 Cities HOSTED Summer Olympics", "Sports"),
    @("16 Cities HOSTED Winter Olympics", "Sports")
)

$startRow = 58
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $rows[$i][0]
    $ws.Range("B$r").Value = $rows[$i][1]
}

# The table's range previously spanned to the sheet's near-bottom row;
# shrink it by one row to match the resulting layout.
$lo = $ws.ListObjects.Item("Table1")
[void]$lo.Resize($ws.Range("A1:B1048567"))

# Update the active selection left after editing.
[void]$ws.Range("D66").Select()
